# [Add] Summon Level process Apply
#  - Fix duplicate-probability issue in the probability table
#  - Add summon levels 3-5 so the drop probability changes as the gacha level changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update probability (C column) for existing rows 3-14 ---
$ws.Range("C3").Value = 9995
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 9995
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1

# --- Step 2: add new rows 15-32 (SummonGrade levels 3,4,5) ---
$ws.Range("A15").Value = 3
$ws.Range("B15").Formula = "=B9"
$ws.Range("C15").Value = 1

$ws.Range("A16").Formula = "=A15"
$ws.Range("B16").Formula = "=B10"
$ws.Range("C16").Value = 1

$ws.Range("A17").Formula = "=A16"
$ws.Range("B17").Formula = "=B11"
$ws.Range("C17").Value = 9995

$ws.Range("A18").Formula = "=A17"
$ws.Range("B18").Formula = "=B12"
$ws.Range("C18").Value = 1

$ws.Range("A19").Formula = "=A18"
$ws.Range("B19").Formula = "=B13"
$ws.Range("C19").Value = 1

$ws.Range("A20").Formula = "=A19"
$ws.Range("B20").Formula = "=B14"
$ws.Range("C20").Value = 1

$ws.Range("A21").Value = 4
$ws.Range("B21").Formula = "=B15"
$ws.Range("C21").Value = 1

$ws.Range("A22").Formula = "=A21"
$ws.Range("B22").Formula = "=B16"
$ws.Range("C22").Value = 1

$ws.Range("A23").Formula = "=A22"
$ws.Range("B23").Formula = "=B17"
$ws.Range("C23").Value = 1

$ws.Range("A24").Formula = "=A23"
$ws.Range("B24").Formula = "=B18"
$ws.Range("C24").Value = 9995

$ws.Range("A25").Formula = "=A24"
$ws.Range("B25").Formula = "=B19"
$ws.Range("C25").Value = 1

$ws.Range("A26").Formula = "=A25"
$ws.Range("B26").Formula = "=B20"
$ws.Range("C26").Value = 1

$ws.Range("A27").Value = 5
$ws.Range("B27").Formula = "=B21"
$ws.Range("C27").Value = 1

$ws.Range("A28").Formula = "=A27"
$ws.Range("B28").Formula = "=B22"
$ws.Range("C28").Value = 1

$ws.Range("A29").Formula = "=A28"
$ws.Range("B29").Formula = "=B23"
$ws.Range("C29").Value = 1

$ws.Range("A30").Formula = "=A29"
$ws.Range("B30").Formula = "=B24"
$ws.Range("C30").Value = 1

$ws.Range("A31").Formula = "=A30"
$ws.Range("B31").Formula = "=B25"
$ws.Range("C31").Value = 9995

$ws.Range("A32").Formula = "=A31"
$ws.Range("B32").Formula = "=B26"
$ws.Range("C32").Value = 1

# --- Step 3: resize the table to cover the new rows (A2:C32) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:C32"))

# --- Step 4: move the active selection (matches the saved workbook state) ---
$ws.Range("G20").Select()
